# Apply the "Updated image and table" edit to the folders table on Sheet1.
#
# The table lists top-level "data" sub-folders. The folder names and their
# descriptions are being refreshed, and the now-unused "docs"/"sandbox" row
# is removed (table shrinks from 6 data rows to 5 data rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update folder names (column B) ---------------------------------------
# Edited in this order so the workbook's shared-string table is rebuilt with
# the same ordering as the authored change.
$ws.Range("B3").Value = "data-lake"
$ws.Range("B2").Value = "data-governance"
$ws.Range("B4").Value = "data-mapping"
$ws.Range("B5").Value = "data-warehouse"

# --- Update descriptions (column C), top to bottom -------------------------
$ws.Range("C2").Value = "Document repository for data governance"
$ws.Range("C3").Value = "Storage of raw, unstructured data"
$ws.Range("C4").Value = "GIS mapping information pertaining to the files stored in the data-lake directory"
$ws.Range("C5").Value = "Storage for curated datasets."

# --- Replace the old "scripts" row with the new "code" row -----------------
$ws.Range("B6").Value = "code"
$ws.Range("C6").Value = "Storage for importation pipeline code as well as secondary data products"

# --- Remove the now-obsolete last row ("sandbox") ---------------------------
$ws.Rows.Item(7).Delete()

# --- Resize column B to fit the new, longer folder names --------------------
$ws.Columns.Item(2).AutoFit() | Out-Null

# --- Move the active selection to reflect where the user left off ----------
$ws.Range("C10").Select() | Out-Null

Write-Output "Folder table updated"
